$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = 44391
$ws.Cells.Item(2, 10).Value = 400
$ws.Cells.Item(2, 11).Value = 15000
$ws.Cells.Item(2, 12).Value = 15000
$ws.Cells.Item(2, 13).Value = 15000
$ws.Cells.Item(2, 16).Value = 833

# Row 3
$ws.Cells.Item(3, 4).Value = 44238
$ws.Cells.Item(3, 9).Value = 'Primera'
$ws.Cells.Item(3, 10).Value = 300
$ws.Cells.Item(3, 11).Value = 12000
$ws.Cells.Item(3, 12).Value = 12000
$ws.Cells.Item(3, 13).Value = 12000
$ws.Cells.Item(3, 16).Value = 667

# Row 4
$ws.Cells.Item(4, 4).Value = 44238
$ws.Cells.Item(4, 9).Value = 'Segunda'
$ws.Cells.Item(4, 10).Value = 200
$ws.Cells.Item(4, 11).Value = 10000
$ws.Cells.Item(4, 12).Value = 10000
$ws.Cells.Item(4, 13).Value = 10000
$ws.Cells.Item(4, 16).Value = 556

# Row 5
$ws.Cells.Item(5, 4).Value = 44238
$ws.Cells.Item(5, 9).Value = 'Tercera'
$ws.Cells.Item(5, 10).Value = 50
$ws.Cells.Item(5, 11).Value = 8000
$ws.Cells.Item(5, 12).Value = 8000
$ws.Cells.Item(5, 13).Value = 8000
$ws.Cells.Item(5, 16).Value = 444

# Row 6
$ws.Cells.Item(6, 4).Value = 44243
$ws.Cells.Item(6, 9).Value = 'Especial'
$ws.Cells.Item(6, 10).Value = 300
$ws.Cells.Item(6, 11).Value = 12000
$ws.Cells.Item(6, 12).Value = 12000
$ws.Cells.Item(6, 13).Value = 12000
$ws.Cells.Item(6, 16).Value = 667

# Row 7
$ws.Cells.Item(7, 4).Value = 44243
$ws.Cells.Item(7, 9).Value = 'Primera'
$ws.Cells.Item(7, 10).Value = 300
$ws.Cells.Item(7, 11).Value = 10000
$ws.Cells.Item(7, 12).Value = 10000
$ws.Cells.Item(7, 13).Value = 10000
$ws.Cells.Item(7, 16).Value = 556

# Row 8
$ws.Cells.Item(8, 4).Value = 44243
$ws.Cells.Item(8, 9).Value = 'Segunda'
$ws.Cells.Item(8, 10).Value = 150
$ws.Cells.Item(8, 11).Value = 8000
$ws.Cells.Item(8, 12).Value = 8000
$ws.Cells.Item(8, 13).Value = 8000
$ws.Cells.Item(8, 16).Value = 444

# Row 9
$ws.Cells.Item(9, 4).Value = 44383
$ws.Cells.Item(9, 11).Value = 16000
$ws.Cells.Item(9, 12).Value = 16000
$ws.Cells.Item(9, 13).Value = 16000
$ws.Cells.Item(9, 14).Value = '$/bandeja 18 kilos'
$ws.Cells.Item(9, 16).Value = 889

# Row 10
$ws.Cells.Item(10, 4).Value = 44383
$ws.Cells.Item(10, 9).Value = 'Segunda'
$ws.Cells.Item(10, 11).Value = 12000
$ws.Cells.Item(10, 12).Value = 12000
$ws.Cells.Item(10, 13).Value = 12000
$ws.Cells.Item(10, 16).Value = 667

# Row 11
$ws.Cells.Item(11, 4).Value = 44229
$ws.Cells.Item(11, 10).Value = 200
$ws.Cells.Item(11, 11).Value = 15000
$ws.Cells.Item(11, 12).Value = 15000
$ws.Cells.Item(11, 13).Value = 15000
$ws.Cells.Item(11, 16).Value = 833

# Row 12
$ws.Cells.Item(12, 4).Value = 44628
$ws.Cells.Item(12, 11).Value = 15000
$ws.Cells.Item(12, 12).Value = 15000
$ws.Cells.Item(12, 13).Value = 15000
$ws.Cells.Item(12, 16).Value = 833

# Row 13
$ws.Cells.Item(13, 4).Value = 44235
$ws.Cells.Item(13, 9).Value = 'Primera'
$ws.Cells.Item(13, 10).Value = 400
$ws.Cells.Item(13, 11).Value = 13000
$ws.Cells.Item(13, 12).Value = 13000
$ws.Cells.Item(13, 13).Value = 13000
$ws.Cells.Item(13, 16).Value = 722

# Row 14
$ws.Cells.Item(14, 4).Value = 44235
$ws.Cells.Item(14, 9).Value = 'Segunda'
$ws.Cells.Item(14, 10).Value = 200
$ws.Cells.Item(14, 11).Value = 11000
$ws.Cells.Item(14, 12).Value = 11000
$ws.Cells.Item(14, 13).Value = 11000
$ws.Cells.Item(14, 16).Value = 611

# Row 15
$ws.Cells.Item(15, 4).Value = 44235
$ws.Cells.Item(15, 9).Value = 'Tercera'
$ws.Cells.Item(15, 10).Value = 100
$ws.Cells.Item(15, 11).Value = 9000
$ws.Cells.Item(15, 12).Value = 9000
$ws.Cells.Item(15, 13).Value = 9000
$ws.Cells.Item(15, 16).Value = 500

# Row 16
$ws.Cells.Item(16, 4).Value = 44631
$ws.Cells.Item(16, 10).Value = 300
$ws.Cells.Item(16, 11).Value = 15000
$ws.Cells.Item(16, 12).Value = 15000
$ws.Cells.Item(16, 13).Value = 15000
$ws.Cells.Item(16, 16).Value = 833

# Row 17
$ws.Cells.Item(17, 4).Value = 44635
$ws.Cells.Item(17, 9).Value = 'Primera'
$ws.Cells.Item(17, 10).Value = 300
$ws.Cells.Item(17, 11).Value = 15000
$ws.Cells.Item(17, 12).Value = 15000
$ws.Cells.Item(17, 13).Value = 15000
$ws.Cells.Item(17, 16).Value = 833

# Row 18
$ws.Cells.Item(18, 4).Value = 44614
$ws.Cells.Item(18, 14).Value = '$/caja 18 kilos granel'

# Row 20
$ws.Cells.Item(20, 4).Value = 44396
$ws.Cells.Item(20, 10).Value = 250
$ws.Cells.Item(20, 11).Value = 15000
$ws.Cells.Item(20, 12).Value = 15000
$ws.Cells.Item(20, 13).Value = 15000
$ws.Cells.Item(20, 16).Value = 833

# Row 21
$ws.Cells.Item(21, 4).Value = 44396
$ws.Cells.Item(21, 10).Value = 150
$ws.Cells.Item(21, 11).Value = 12000
$ws.Cells.Item(21, 12).Value = 12000
$ws.Cells.Item(21, 13).Value = 12000
$ws.Cells.Item(21, 16).Value = 667

# Row 22
$ws.Cells.Item(22, 4).Value = 44596
$ws.Cells.Item(22, 9).Value = 'Primera'
$ws.Cells.Item(22, 10).Value = 150
$ws.Cells.Item(22, 11).Value = 14000
$ws.Cells.Item(22, 12).Value = 14000
$ws.Cells.Item(22, 13).Value = 14000
$ws.Cells.Item(22, 16).Value = 778

# Row 23
$ws.Cells.Item(23, 4).Value = 44630

# Row 25
$ws.Cells.Item(25, 4).Value = 44245
$ws.Cells.Item(25, 9).Value = 'Primera'

# Row 26
$ws.Cells.Item(26, 4).Value = 44245
$ws.Cells.Item(26, 9).Value = 'Segunda'
$ws.Cells.Item(26, 10).Value = 200

# Row 27
$ws.Cells.Item(27, 4).Value = 44249
$ws.Cells.Item(27, 9).Value = 'Primera'
$ws.Cells.Item(27, 10).Value = 400
$ws.Cells.Item(27, 11).Value = 12000
$ws.Cells.Item(27, 12).Value = 12000
$ws.Cells.Item(27, 13).Value = 12000
$ws.Cells.Item(27, 16).Value = 667

# Row 28
$ws.Cells.Item(28, 4).Value = 44249
$ws.Cells.Item(28, 9).Value = 'Segunda'
$ws.Cells.Item(28, 10).Value = 200
$ws.Cells.Item(28, 11).Value = 10000
$ws.Cells.Item(28, 12).Value = 10000
$ws.Cells.Item(28, 13).Value = 10000
$ws.Cells.Item(28, 16).Value = 556
